$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 59: homework entry added 2020-02-13 (commit: "Added homework 13Feb 2020")
$row = 59

# Columns B (date-like) and C (zero-padded id) must stay text, not be
# auto-coerced into a date serial / number by Excel's smart input parsing.
# Force a text number format before writing, then restore the default
# "Normal" style so no stray formatting is left on the cell.
$ws.Range("B$row" + ":C$row").NumberFormat = "@"

$ws.Range("A$row").Value = 1581552000
$ws.Range("B$row").Value = "2020-02-13"
$ws.Range("C$row").Value = "0216"
$ws.Range("D$row").Value = "SPRING"
$ws.Range("E$row").Value = 0.225
$ws.Range("F$row").Value = 0.225
$ws.Range("G$row").Value = 0.22
$ws.Range("H$row").Value = 0.22
$ws.Range("I$row").Value = 890600

$ws.Range("B$row" + ":C$row").Style = "Normal"
